# Updated cryptos list on Mon Oct 21 09:46:18 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for existing rows and
# inserts a new "Binance-PegBSC-USD" entry at row 30, shifting the rest
# of the ranking table (rows 30-51) down by one (Cronos drops off the
# bottom of the A1:E51 range).
#
# Numeric-looking Price strings (e.g. "605.97") are written with a
# temporary Text number format so Excel stores them as strings (matching
# the source data's inlineStr type) instead of auto-coercing them to
# numbers; the style is then reset to "Normal" so no stray formatting is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.279.29"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.708.34"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.42%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").Value = "2.707.17"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  +3.59%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "3.220.01"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "68.173.99"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "2.705.37"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "370.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("D28").Value = "2.844.19"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "575.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.376"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "0.0₆0310"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.593"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "154.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.51%  "
